# Refresh the cryptocurrency price/volume snapshot (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells that look like plain decimals (e.g. "19.94") must stay TEXT,
# matching the source workbook. A leading single-quote is the standard Excel
# "treat this as text" marker; in a single-quoted PS string a doubled '' is
# how you embed one literal apostrophe, so '' + value => a text-marked cell.

$ws.Range("D2").Value = '27.176.85'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '1.645.97'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '''217.31'
$ws.Range("E5").Value = '  -1.39%  '
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("D10").Value = '''19.94'
$ws.Range("E10").Value = '  -0.25%  '
$ws.Range("D11").Value = '''0.0844'
$ws.Range("E11").Value = '  -0.84%  '
$ws.Range("D12").Value = '1.875.08'
$ws.Range("E12").Value = '  -1.04%  '
$ws.Range("D13").Value = '1.665.41'
$ws.Range("E13").Value = '  +0.01%  '
$ws.Range("E14").Value = '  -2.51%  '
$ws.Range("D15").Value = '''0.541'
$ws.Range("E15").Value = '  +0.88%  '
$ws.Range("D16").Value = '''67.40'
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("D17").Value = '27.139.13'
$ws.Range("E17").Value = '  -0.70%  '
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("D19").Value = '''217.90'
$ws.Range("E19").Value = '  -2.50%  '
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").Value = '''6.83'
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("E22").Value = '  -0.56%  '
$ws.Range("E23").Value = '  +0.48%  '
$ws.Range("D24").Value = '''9.18'
$ws.Range("E24").Value = '  -1.32%  '
$ws.Range("D25").Value = '''147.45'
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").Value = '''7.44'
$ws.Range("E27").Value = '  -0.28%  '
$ws.Range("E28").Value = '  -1.33%  '
$ws.Range("D29").Value = '''15.73'
$ws.Range("E29").Value = '  -2.29%  '
$ws.Range("D30").Value = '''0.0505'
$ws.Range("E30").Value = '  -2.19%  '
$ws.Range("E31").Value = '  -1.54%  '
$ws.Range("E32").Value = '  -0.93%  '
$ws.Range("E33").Value = '  +0.47%  '
$ws.Range("D34").Value = '''1.58'
$ws.Range("E34").Value = '  +0.82%  '
$ws.Range("D35").Value = '1.266.59'
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").Value = '''0.0177'
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("E38").Value = '  +0.30%  '
$ws.Range("E39").Value = '  +0.60%  '
$ws.Range("E40").Value = '  -0.20%  '
$ws.Range("D41").Value = '''0.809'
$ws.Range("E41").Value = '  -0.75%  '
$ws.Range("E42").Value = '  +4.48%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").Value = '1.785.13'
$ws.Range("D45").Value = '''62.40'
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("D46").Value = '''91.72'
$ws.Range("E46").Value = '  -0.90%  '
$ws.Range("D47").Value = '''1.61'
$ws.Range("E47").Value = '  -1.22%  '
$ws.Range("D48").Value = '0.0₆0107'
$ws.Range("E48").Value = '  +15.72%  '
$ws.Range("E49").Value = '  -0.96%  '
$ws.Range("D50").Value = '''7.68'
$ws.Range("E50").Value = '  -0.16%  '
$ws.Range("D51").Value = '''0.0972'
$ws.Range("E51").Value = '  -1.40%  '
